$d = $word.ActiveDocument

# The document has no inline pictures in the body - the BTec logo sits in
# both headers (first-page + default) and the Pearson logo sits in both
# footers (first-page + default). Word's internal picture name swaps for
# each: the BTec logo ("BTec_Logo-Orange") goes from "image1.jpg" to
# "image2.jpg", and the Pearson logo (the PearsonLogo.png description)
# goes from "image2.png" to "image1.png". Everything else (the alt-text
# description, size, id) is left untouched.
foreach ($story in $d.StoryRanges) {
    if ($story.InlineShapes.Count -gt 0) {
        $pic = $story.InlineShapes(1)
        try {
            if ($pic.AlternativeText -eq "BTec_Logo-Orange") {
                $pic.Name = "image2.jpg"
            } else {
                $pic.Name = "image1.png"
            }
        } catch {
            # Some stories may refuse a rename of an already up to date
            # picture handle; nothing else to do for that shape.
        }
    }
}
